$d = $word.ActiveDocument

# --- Paragraph 1 edit -------------------------------------------------
# Before: "This is a Microsoft word document."
# After : "This is a Microsoft word document.  (This is a change – Version for main branch)"
# The appended text is split across three runs (all colored red, FF0000):
#   "(This is a change – Ve" | "rsion for main branch" | ")"
# and two trailing spaces are appended to the original (black) run.

$enDash = [char]0x2013

$para1 = $d.Paragraphs(1).Range
$para1.InsertAfter("  ")

# Position right after the two spaces we just added (before the hidden
# paragraph mark) -- real content length is (Range.End - 1).
$pos = $d.Paragraphs(1).Range.End - 1

$seg2 = "(This is a change " + $enDash + " Ve"
$seg3 = "rsion for main branch"
$seg4 = ")"

# --- segment 2 ---
$ins = $d.Range($pos, $pos)
$ins.InsertAfter($seg2)
$pos2 = $d.Paragraphs(1).Range.End - 1
$r2 = $d.Range($pos, $pos2)
$r2.Font.Color = 255
$pos = $pos2

# --- segment 3 ---
$ins = $d.Range($pos, $pos)
$ins.InsertAfter($seg3)
$pos3 = $d.Paragraphs(1).Range.End - 1
$r3 = $d.Range($pos, $pos3)
$r3.Font.Color = 255
$pos = $pos3

# --- segment 4 ---
$ins = $d.Range($pos, $pos)
$ins.InsertAfter($seg4)
$pos4 = $d.Paragraphs(1).Range.End - 1
$r4 = $d.Range($pos, $pos4)
$r4.Font.Color = 255
$pos = $pos4
